$d = $word.ActiveDocument

# 1. Insert two new achievement bullets before the "Developed automatic build..."
#    bullet in the Amino Communications section.
$r = $d.Content
$r.Find.Execute("Developed automatic build and release systems", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$para.Range.InsertBefore("Implemented rsyslog support to the AminoMOVE CloudTV media backend, enabling its output to be fed to the log collection and analysis tools Splunk and Unomaly allowing production deployment issues to be field issues to be rapidly identified and understood`rDelivered a complete set of demonstrator services used to support the sales team in winning opportunities; including developing the ground-breaking concept of using an Amino STB as the cable to IP turn-around acting as the content source for the Amino's first AWS hosted CloudTV service`r")

# 2. Remove the "Platform ports Amino's AmiNET..." bullet that followed the
#    "Developed the company's first Scrum tool..." bullet.
$r2 = $d.Content
$r2.Find.Execute("Platform ports Amino", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $r2.Paragraphs(1)
$para2.Range.Delete()

# 3. Tweak the Global Graphics company summary wording.
$d.Content.Find.Execute("PDF workflow software for Windows", $true, $false, $false, $false, $false, $true, 1, $false, "PDF workflow software using for Windows", 2)

# 4. Insert two new achievement bullets after the Global Graphics company
#    summary paragraph. A new paragraph mark is split off the end of the
#    company_summary paragraph (so it starts out empty, styled
#    company_summary), then it is restyled to achievement_bullet and
#    filled with both new bullets before a final leftover empty
#    paragraph is removed.
$r3 = $d.Content
$r3.Find.Execute("Electronic document software; Managed team to deliver PDF workflow software using for Windows", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos3 = $r3.End
$d.Range($endPos3, $endPos3).InsertParagraphAfter()

$newPara3 = $d.Range($endPos3 + 1, $endPos3 + 1).Paragraphs(1)
$newPara3.Style = "achievement_bullet"
$newPara3.Range.InsertBefore("Formulated framework for developer tests, using QT’s built-in javascript technology. Used to implement tests during module development which were then built into a suite for avoiding regressions in subsequent iterations`rDefined change request tracking processes and then evolved it throughout the project to support team’s constant process improvement`r")

$r3c = $d.Content
$r3c.Find.Execute("Defined change request tracking processes and then evolved it throughout the project to support team’s constant process improvement", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3c.Paragraphs(1).Next().Range.Delete()

# 5. Remove the "Used SQL and Perl skills..." bullet.
$r4 = $d.Content
$r4.Find.Execute("Used SQL and Perl skills to write tool to export issue tracking history", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para4 = $r4.Paragraphs(1)
$para4.Range.Delete()

# 6. Remove the "Recruited test team leader..." bullet.
$r5 = $d.Content
$r5.Find.Execute("Recruited test team leader, software and test engineers", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para5 = $r5.Paragraphs(1)
$para5.Range.Delete()
